$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")
$ws.Activate()

# Insert a new row above row 12 (shifts existing rows 12+ down by one),
# making room for the new "keywordDocumentPath" constant.
$null = $ws.Rows.Item(12).Insert()

$ws.Range("A12").Value = "keywordDocumentPath"
$ws.Range("B12").Value = "DocumentProcessing\keyword.json"

$null = $ws.Range("B11").Select()
